# "Duplicate Detection" workbook edit
#
# The author duplicated the "Sheet1" Power Query (Excel creates a query
# named "Sheet1 (2)" together with a new connection/query table) and the
# duplicate's output landed on a brand-new worksheet. That new worksheet
# was inserted right before the previously-active "Sheet2" tab (which is
# why Excel names it "Sheet3" even though it is the 2nd tab), becomes the
# new active tab, and receives the query's "ExternalData_2" external-data
# range name anchored at A1. "Sheet2" keeps its own "ExternalData_1" name,
# whose localSheetId shifts because of the newly inserted sheet.

$wb = $excel.ActiveWorkbook

# Insert a new worksheet. With no target sheet specified, Excel inserts it
# immediately before the active sheet ("Sheet2" here), and it becomes the
# new active sheet itself - exactly matching Sheet1, Sheet3, Sheet2 order.
$ws3 = $wb.Worksheets.Add()

# The query output is a single cell holding the query/table name.
$ws3.Range("A1").Value = "ExternalData_2"

# Excel records the query's output range as a sheet-scoped defined name
# named after the query, e.g. "ExternalData_2" -> Sheet3!$A$1.
$ws3.Names.Add("ExternalData_2", "=" + $ws3.Name + "!`$A`$1")
